$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-21: news items changed (title, relative time, url, relevance, keyword)
$ws.Cells.Item(2, 1).Value = "Lalin di Terminal Kampung Rambutan Macet Imbas Proyek Galian"
$ws.Cells.Item(2, 2).Value = "6 menit yang lalu"
$ws.Cells.Item(2, 4).Value = "https://news.detik.com/berita/d-8135129/lalin-di-terminal-kampung-rambutan-macet-imbas-proyek-galian"
$ws.Cells.Item(2, 5).Value = $false
$ws.Cells.Item(2, 6).Value = ""

$ws.Cells.Item(3, 1).Value = "Bareskrim Tangkap Kurir Narkoba di Jakut, 10 Bungkus Sabu Disita"
$ws.Cells.Item(3, 2).Value = "16 menit yang lalu"
$ws.Cells.Item(3, 4).Value = "https://news.detik.com/berita/d-8135121/bareskrim-tangkap-kurir-narkoba-di-jakut-10-bungkus-sabu-disita"
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(3, 6).Value = ""

$ws.Cells.Item(4, 1).Value = "Cerita Jaksa di Bantaeng Alami Intimidasi Saat Tangani Perkara Besar"
$ws.Cells.Item(4, 2).Value = "43 menit yang lalu"
$ws.Cells.Item(4, 4).Value = "https://news.detik.com/berita/d-8135073/cerita-jaksa-di-bantaeng-alami-intimidasi-saat-tangani-perkara-besar"
$ws.Cells.Item(4, 5).Value = $false
$ws.Cells.Item(4, 6).Value = ""

$ws.Cells.Item(5, 1).Value = "Kementrans Akan Perluas Program Pengiriman Transmigran RI Ke Jepang"
$ws.Cells.Item(5, 2).Value = "48 menit yang lalu"
$ws.Cells.Item(5, 4).Value = "https://news.detik.com/berita/d-8135087/kementrans-akan-perluas-program-pengiriman-transmigran-ri-ke-jepang"
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = ""

$ws.Cells.Item(6, 1).Value = "Legislator Setuju Koki Makan Bergizi Gratis Harus Punya Pengalaman"
$ws.Cells.Item(6, 2).Value = "50 menit yang lalu"
$ws.Cells.Item(6, 4).Value = "https://news.detik.com/berita/d-8135085/legislator-setuju-koki-makan-bergizi-gratis-harus-punya-pengalaman"
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(6, 6).Value = ""

$ws.Cells.Item(7, 1).Value = "Upaya Kejaksaan Edukasi Restorative Justice ke Warga Jeneponto Sulsel"
$ws.Cells.Item(7, 2).Value = "58 menit yang lalu"
$ws.Cells.Item(7, 4).Value = "https://news.detik.com/berita/d-8134826/upaya-kejaksaan-edukasi-restorative-justice-ke-warga-jeneponto-sulsel"
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = ""

$ws.Cells.Item(8, 1).Value = "Dapur-Chef Dievaluasi Buntut Menu MBG Bikin Siswa Keracunan"
$ws.Cells.Item(8, 2).Value = "1 jam yang lalu"
$ws.Cells.Item(8, 4).Value = "https://news.detik.com/berita/d-8135062/dapur-chef-dievaluasi-buntut-menu-mbg-bikin-siswa-keracunan"
$ws.Cells.Item(8, 5).Value = $true
$ws.Cells.Item(8, 6).Value = "mbg"

$ws.Cells.Item(9, 1).Value = "Juru Masak MBG Dievaluasi"
$ws.Cells.Item(9, 2).Value = "1 jam yang lalu"
$ws.Cells.Item(9, 4).Value = "https://news.detik.com/berita/d-8134900/juru-masak-mbg-dievaluasi"
$ws.Cells.Item(9, 5).Value = $true
$ws.Cells.Item(9, 6).Value = "mbg"

$ws.Cells.Item(10, 1).Value = "Lalin di Sejumlah Ruas Tol Arah Jakarta Padat, Ini Titiknya"
$ws.Cells.Item(10, 2).Value = "1 jam yang lalu"
$ws.Cells.Item(10, 4).Value = "https://news.detik.com/berita/d-8135054/lalin-di-sejumlah-ruas-tol-arah-jakarta-padat-ini-titiknya"
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = ""

$ws.Cells.Item(11, 1).Value = "Sejumlah Gerbang Tol Dalam Kota Dibuka Parsial Pagi Ini, Lalin Padat"
$ws.Cells.Item(11, 2).Value = "1 jam yang lalu"
$ws.Cells.Item(11, 4).Value = "https://news.detik.com/berita/d-8135051/sejumlah-gerbang-tol-dalam-kota-dibuka-parsial-pagi-ini-lalin-padat"
$ws.Cells.Item(11, 5).Value = $false
$ws.Cells.Item(11, 6).Value = ""

$ws.Cells.Item(12, 1).Value = "Update Terbaru Kasus Anggota TNI Pukul Staf Zaskia Adya Mecca"
$ws.Cells.Item(12, 2).Value = "1 jam yang lalu"
$ws.Cells.Item(12, 4).Value = "https://news.detik.com/berita/d-8135048/update-terbaru-kasus-anggota-tni-pukul-staf-zaskia-adya-mecca"
$ws.Cells.Item(12, 5).Value = $false
$ws.Cells.Item(12, 6).Value = ""

$ws.Cells.Item(13, 1).Value = "JPPI Nilai Sertifikat Higienis dan CCTV Belum Bisa Selesaikan Persoalan MBG"
$ws.Cells.Item(13, 2).Value = "2 jam yang lalu"
$ws.Cells.Item(13, 4).Value = "https://news.detik.com/berita/d-8135046/jppi-nilai-sertifikat-higienis-dan-cctv-belum-bisa-selesaikan-persoalan-mbg"
$ws.Cells.Item(13, 5).Value = $true
$ws.Cells.Item(13, 6).Value = "mbg"

$ws.Cells.Item(14, 1).Value = "Terungkap Bakteri 'Biang Kerok' Keracunan MBG di Bandung Barat"
$ws.Cells.Item(14, 2).Value = "2 jam yang lalu"
$ws.Cells.Item(14, 4).Value = "https://news.detik.com/berita/d-8135037/terungkap-bakteri-biang-kerok-keracunan-mbg-di-bandung-barat"
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = "mbg"

$ws.Cells.Item(15, 1).Value = "Hari Jantung Sedunia 2025: Tema, Tujuan, dan Cara Merayakan"
$ws.Cells.Item(15, 2).Value = "2 jam yang lalu"
$ws.Cells.Item(15, 4).Value = "https://news.detik.com/berita/d-8134113/hari-jantung-sedunia-2025-tema-tujuan-dan-cara-merayakan"
$ws.Cells.Item(15, 5).Value = $false
$ws.Cells.Item(15, 6).Value = ""

$ws.Cells.Item(16, 1).Value = "Anggota DPR Minta Penanganan Keracunan MBG Tak Cuma CCTV-Sertifikat Higienis"
$ws.Cells.Item(16, 2).Value = "2 jam yang lalu"
$ws.Cells.Item(16, 4).Value = "https://news.detik.com/berita/d-8135035/anggota-dpr-minta-penanganan-keracunan-mbg-tak-cuma-cctv-sertifikat-higienis"
$ws.Cells.Item(16, 5).Value = $true
$ws.Cells.Item(16, 6).Value = "mbg"

$ws.Cells.Item(17, 1).Value = "4 Instruksi Prabowo Saat Panggil Kepala BGN Usai Marak Keracunan MBG"
$ws.Cells.Item(17, 2).Value = "2 jam yang lalu"
$ws.Cells.Item(17, 4).Value = "https://news.detik.com/berita/d-8135031/4-instruksi-prabowo-saat-panggil-kepala-bgn-usai-marak-keracunan-mbg"
$ws.Cells.Item(17, 5).Value = $true
$ws.Cells.Item(17, 6).Value = "mbg"

$ws.Cells.Item(18, 1).Value = "29 September Hari Sarjana Nasional, Ini Sejarah dan Cara Merayakannya"
$ws.Cells.Item(18, 2).Value = "3 jam yang lalu"
$ws.Cells.Item(18, 4).Value = "https://news.detik.com/berita/d-8134101/29-september-hari-sarjana-nasional-ini-sejarah-dan-cara-merayakannya"
$ws.Cells.Item(18, 5).Value = $false
$ws.Cells.Item(18, 6).Value = ""

$ws.Cells.Item(19, 1).Value = "CISDI Dorong Pemerintah Susun Perpres Perbaikan Menyeluruh MBG"
$ws.Cells.Item(19, 2).Value = "3 jam yang lalu"
$ws.Cells.Item(19, 4).Value = "https://news.detik.com/berita/d-8135020/cisdi-dorong-pemerintah-susun-perpres-perbaikan-menyeluruh-mbg"
$ws.Cells.Item(19, 5).Value = $true
$ws.Cells.Item(19, 6).Value = "mbg"

$ws.Cells.Item(20, 1).Value = "400 Rumah Hangus Usai Kebakaran Landa Permukiman di Tamansari Jakbar"
$ws.Cells.Item(20, 2).Value = "4 jam yang lalu"
$ws.Cells.Item(20, 4).Value = "https://news.detik.com/berita/d-8135009/400-rumah-hangus-usai-kebakaran-landa-permukiman-di-tamansari-jakbar"
$ws.Cells.Item(20, 5).Value = $false
$ws.Cells.Item(20, 6).Value = ""

$ws.Cells.Item(21, 1).Value = "Cak Imin Sebut MBG Jadi Solusi Agar Anak Tak Banyak Konsumsi MSG"
$ws.Cells.Item(21, 2).Value = "4 jam yang lalu"
$ws.Cells.Item(21, 4).Value = "https://news.detik.com/berita/d-8135005/cak-imin-sebut-mbg-jadi-solusi-agar-anak-tak-banyak-konsumsi-msg"
$ws.Cells.Item(21, 5).Value = $true
$ws.Cells.Item(21, 6).Value = "mbg"

# Rows 22-142: menu/footer links - only the scrape date (column B) moved forward one day
$dateRange = $ws.Range("B22:B142")
$dateRange.NumberFormat = "@"
for ($r = 22; $r -le 142; $r++) {
    $ws.Cells.Item($r, 2).Value = "2025-09-29"
}

Write-Host "Edit complete"
